$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New loading-percent values for rows 2-25 (columns B,C,D,F,G,H,K,O)
$data = @{
    2 = @{ 2=8.997321071046972; 3=6.919270698895518; 4=4.700355325041238; 6=20.0320782451536; 7=21.35332168543748; 8=12.64010216132666; 11=9.026342301698193; 15=18.21877007728061 }
    3 = @{ 2=8.593600640786612; 3=6.80365787432458; 4=4.577245143934905; 6=20.14374782432396; 7=21.55010247829403; 8=12.69944213428341; 11=8.71763340129986; 15=18.33375722584825 }
    4 = @{ 2=8.335900059793278; 3=6.731836500228132; 4=4.499214202841038; 6=20.21839751029265; 7=21.67942905547092; 8=12.73793592786366; 11=8.521160768418769; 15=18.40859208710413 }
    5 = @{ 2=8.228546149228047; 3=6.702388778868662; 4=4.466831102332942; 6=20.25034303777106; 7=21.73426095189424; 8=12.75414105733104; 11=8.439438627085172; 15=18.44015264227836 }
    6 = @{ 2=8.210583000873999; 3=6.69748902957558; 4=4.461419476165821; 6=20.25573954639288; 7=21.74349418855765; 8=12.75686325510321; 11=8.425771196657029; 15=18.4454575719236 }
    7 = @{ 2=8.334461532759885; 3=6.731440046371118; 4=4.498779801327661; 6=20.21882216972956; 7=21.68015992224237; 8=12.73815237493625; 11=8.520065231613986; 15=18.40901341198292 }
    8 = @{ 2=8.860227065234618; 3=6.879598348506219; 4=4.658430184744993; 6=20.06931632997402; 7=21.41940232279584; 8=12.6601358144924; 11=8.921383763656241; 15=18.25753952792454 }
    9 = @{ 2=9.808602843757157; 3=7.16224123440043; 4=4.950979057601574; 6=19.8246199236234; 7=20.97587498347187; 8=12.52344045089448; 11=9.650212190735271; 15=17.99405832655994 }
    10 = @{ 2=10.44970025933289; 3=7.363451649382577; 4=5.15200380099698; 6=19.67468210636081; 7=20.69186013613089; 8=12.43288674866353; 11=10.14649417033902; 15=17.82090552800811 }
    11 = @{ 2=10.72841699007262; 3=7.453263881153285; 4=5.240180567707514; 6=19.61301099173721; 7=20.57185475723108; 8=12.39382412562459; 11=10.36311567891324; 15=17.74656393657688 }
    12 = @{ 2=10.83204590228777; 3=7.487000834429337; 4=5.273081927013843; 6=19.59060258403747; 7=20.52774414238996; 8=12.37933767888958; 11=10.44378559217916; 15=17.71904909943488 }
    13 = @{ 2=10.80981354123406; 3=7.479747514594052; 4=5.266018073505936; 6=19.59538652004073; 7=20.53718470145091; 8=12.38244400786166; 11=10.42647303184873; 15=17.72494659027782 }
    14 = @{ 2=10.73698127626202; 3=7.456045054334637; 4=5.242897290145887; 6=19.61114846781032; 7=20.56819897637172; 8=12.39262619360995; 11=10.36977993730727; 15=17.74428751240414 }
    15 = @{ 2=10.69211848573564; 3=7.44149032091624; 4=5.228670893873173; 6=19.62092634909002; 7=20.58737001021137; 8=12.39890287065901; 11=10.33487541973538; 15=17.75621731096494 }
    16 = @{ 2=10.43122035917527; 3=7.357545378719377; 4=5.146173766905724; 6=19.67884443184051; 7=20.69988868924606; 8=12.4354824002988; 11=10.13214949024998; 15=17.82585301228835 }
    17 = @{ 2=10.26781461682259; 3=7.305588753318927; 4=5.094712712294606; 6=19.71605312909458; 7=20.77127764570071; 8=12.45846799763298; 11=10.00540838890021; 15=17.86970609758902 }
    18 = @{ 2=10.17261444211378; 3=7.27554429421557; 4=5.064806861511704; 6=19.73806951420103; 7=20.81320344573104; 8=12.47188929670997; 11=9.931652320882311; 15=17.89534585821732 }
    19 = @{ 2=10.14017477146598; 3=7.265345017282389; 4=5.054629168365179; 6=19.74562933890371; 7=20.82754700511228; 8=12.4764679905852; 11=9.906533922441785; 15=17.90409859205229 }
    20 = @{ 2=10.2853354740983; 3=7.311136414536767; 4=5.100222726932161; 6=19.71202851935981; 7=20.76358859374021; 8=12.45600038603565; 11=10.01898933875852; 15=17.86499474147202 }
    21 = @{ 2=10.75842625906646; 3=7.463014648736296; 4=5.249701845365147; 6=19.60649311085651; 7=20.55905307295118; 8=12.38962714691696; 11=10.38646931684555; 15=17.73858933274785 }
    22 = @{ 2=11.05643598695265; 3=7.56067368865153; 4=5.344535745377939; 6=19.5430307643392; 7=20.43315116719076; 8=12.34803007614222; 11=10.61869476005571; 15=17.65968746588203 }
    23 = @{ 2=10.89842219997662; 3=7.508705970501435; 4=5.294188534447495; 6=19.57639585034532; 7=20.49963245638942; 8=12.3700684087994; 11=10.49549213251572; 15=17.70145920448343 }
    24 = @{ 2=10.27741820238464; 3=7.308628856606841; 4=5.097732645226571; 6=19.7138461002245; 7=20.7670620634235; 8=12.45711534972913; 11=10.01285216283875; 15=17.86712341228419 }
    25 = @{ 2=9.561491735250826; 3=7.086797060417347; 4=4.874186759183072; 6=19.88559679010333; 7=21.0885500021882; 8=12.55868169954606; 11=9.45968613354388; 15=18.06174780040752 }
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    foreach ($c in $rowVals.Keys) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c]
    }
}
